$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.745.41"
$ws.Range("E2").Value = "  +3.44%  "
$ws.Range("D3").Value = "4.032.29"
$ws.Range("E3").Value = "  +3.29%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'516.46"
$ws.Range("E5").Value = "  -1.84%  "
$ws.Range("D6").Value = "'147.18"
$ws.Range("E6").Value = "  +2.03%  "
$ws.Range("D7").Value = "'0.619"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("E8").Value = "  +0.20%  "
$ws.Range("D9").Value = "'0.733"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "'0.174"
$ws.Range("E10").Value = "  +2.01%  "
$ws.Range("D11").Value = "'0.0000333"
$ws.Range("E11").Value = "  +1.28%  "
$ws.Range("D12").Value = "'46.60"
$ws.Range("E12").Value = "  +10.56%  "
$ws.Range("D13").Value = "'10.71"
$ws.Range("E13").Value = "  +4.86%  "
$ws.Range("D14").Value = "4.672.14"
$ws.Range("E14").Value = "  +3.28%  "
$ws.Range("D15").Value = "4.037.66"
$ws.Range("E15").Value = "  +3.45%  "
$ws.Range("D16").Value = "'21.16"
$ws.Range("E16").Value = "  +7.11%  "
$ws.Range("D17").Value = "'14.21"
$ws.Range("E17").Value = "  +2.60%  "
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").Value = "'0.132"
$ws.Range("E19").Value = "  -2.21%  "
$ws.Range("D20").Value = "71.703.42"
$ws.Range("E20").Value = "  +3.53%  "
$ws.Range("D21").Value = "'437.51"
$ws.Range("E21").Value = "  +1.11%  "
$ws.Range("D22").Value = "'95.96"
$ws.Range("E22").Value = "  +9.59%  "
$ws.Range("D23").Value = "'3.51"
$ws.Range("E23").Value = "  +5.47%  "
$ws.Range("D24").Value = "'14.46"
$ws.Range("E24").Value = "  +1.38%  "
$ws.Range("E25").Value = "  +3.61%  "
$ws.Range("D26").Value = "'4.02"
$ws.Range("E26").Value = "  -0.39%  "
$ws.Range("D27").Value = "'11.17"
$ws.Range("E27").Value = "  +4.99%  "
$ws.Range("D28").Value = "'36.81"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("D29").Value = "'3.08"
$ws.Range("E29").Value = "  +9.55%  "
$ws.Range("D30").Value = "'702.50"
$ws.Range("E30").Value = "  +1.16%  "
$ws.Range("D31").Value = "'13.41"
$ws.Range("E31").Value = "  +1.82%  "

$ws.Range("B32").Value = "NEARProtocol"
$ws.Range("C32").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D32").Value = "'7.19"
$ws.Range("E32").Value = "  +21.96%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.129"
$ws.Range("E33").Value = "  +2.75%  "

$ws.Range("D34").Value = "'67.86"
$ws.Range("E34").Value = "  +0.85%  "
$ws.Range("D35").Value = "0.0₃0896"
$ws.Range("E35").Value = "  +7.74%  "
$ws.Range("D36").Value = "'3.77"
$ws.Range("E36").Value = "  +28.05%  "
$ws.Range("D37").Value = "'0.437"
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("D38").Value = "'40.38"
$ws.Range("E38").Value = "  +0.99%  "
$ws.Range("D39").Value = "'0.154"
$ws.Range("E39").Value = "  +3.53%  "
$ws.Range("E40").Value = "  +0.08%  "
$ws.Range("D41").Value = "'0.997"
$ws.Range("E41").Value = "  -0.32%  "
$ws.Range("D42").Value = "'0.0484"
$ws.Range("E42").Value = "  +0.97%  "
$ws.Range("D43").Value = "'3.16"
$ws.Range("E43").Value = "  +3.11%  "
$ws.Range("D44").Value = "'2.77"
$ws.Range("E44").Value = "  +0.43%  "
$ws.Range("D45").Value = "'3.52"
$ws.Range("E45").Value = "  +4.79%  "
$ws.Range("E46").Value = "  +2.79%  "
$ws.Range("D47").Value = "'3.15"
$ws.Range("E47").Value = "  +3.25%  "
$ws.Range("D48").Value = "'0.000273"
$ws.Range("E48").Value = "  +21.16%  "
$ws.Range("D49").Value = "'9.00"
$ws.Range("E49").Value = "  +6.15%  "
$ws.Range("D50").Value = "'3.31"
$ws.Range("E50").Value = "  +0.88%  "
$ws.Range("D51").Value = "0.0₆0340"
$ws.Range("E51").Value = "  +1.70%  "
